$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1401.091
$ws.Range("J32").Value = 1322
$ws.Range("L32").Value = 1322
$ws.Range("N32").Value = -1974

$ws.Range("H121").Value = 4092
$ws.Range("J121").Value = 4092
$ws.Range("L121").Value = 12276
$ws.Range("N121").Value = -15770

$ws.Range("H132").Value = 105443.31
$ws.Range("I132").Value = 229716.52
$ws.Range("K132").Value = 689149.5599999999
$ws.Range("M132").Value = -686619.5599999999

$ws.Range("H137").Value = 2412.9333
$ws.Range("I137").Value = 1563.091
$ws.Range("J137").Value = 4750
$ws.Range("K137").Value = 4689.272999999999
$ws.Range("L137").Value = 14250
$ws.Range("M137").Value = -2139.272999999999
$ws.Range("N137").Value = -19350

$ws.Range("H140").Value = 52043
$ws.Range("J140").Value = 50047.777
$ws.Range("L140").Value = 50047.777
$ws.Range("N140").Value = -60407.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2744.5
$ws.Range("I32").Value = 2164.5833
$ws.Range("K32").Value = 2164.5833
$ws.Range("M32").Value = -1877.5833

$ws.Range("H43").Value = 8944
$ws.Range("J43").Value = 9758.666999999999
$ws.Range("L43").Value = 9758.666999999999
$ws.Range("N43").Value = -10384.667

$ws.Range("H45").Value = 2774.75
$ws.Range("I45").Value = 1840.1177
$ws.Range("K45").Value = 1840.1177
$ws.Range("M45").Value = -1463.1177

$ws.Range("H74").Value = 4810524.5
$ws.Range("I74").Value = 10001090
$ws.Range("J74").Value = 4445.963
$ws.Range("K74").Value = 10001090
$ws.Range("L74").Value = 4445.963
$ws.Range("M74").Value = -10000216
$ws.Range("N74").Value = -6193.963

$ws.Range("H77").Value = 4810524.5
$ws.Range("I77").Value = 10001090
$ws.Range("J77").Value = 4445.963
$ws.Range("K77").Value = 50005450
$ws.Range("L77").Value = 22229.815
$ws.Range("M77").Value = -50001082
$ws.Range("N77").Value = -30965.815

$ws.Range("H122").Value = 5749.2354
$ws.Range("I122").Value = 3300.125
$ws.Range("K122").Value = 9900.375
$ws.Range("M122").Value = -7450.375

$ws.Range("H132").Value = 3727.2666
$ws.Range("I132").Value = 2981.8086
$ws.Range("J132").Value = 4978.5713
$ws.Range("K132").Value = 8945.425799999999
$ws.Range("L132").Value = 14935.7139
$ws.Range("M132").Value = -6415.425799999999
$ws.Range("N132").Value = -19995.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 27999
$ws.Range("J88").Value = 27999
$ws.Range("L88").Value = 27999
$ws.Range("N88").Value = -28811

$ws.Range("H91").Value = 27999
$ws.Range("J91").Value = 27999
$ws.Range("L91").Value = 27999
$ws.Range("N91").Value = -30807

$ws.Range("H99").Value = 26868090
$ws.Range("I99").Value = 78527880
$ws.Range("K99").Value = 78527880
$ws.Range("M99").Value = -78526382

$ws.Range("H105").Value = 166668350
$ws.Range("I105").Value = 166668350
$ws.Range("K105").Value = 166668350
$ws.Range("M105").Value = -166666603

$ws.Range("H107").Value = 1426.5
$ws.Range("I107").Value = 1046.1177
$ws.Range("K107").Value = 1046.1177
$ws.Range("M107").Value = 873.8823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 55556732
$ws.Range("K31").Value = 55556732
$ws.Range("M31").Value = -55556437

$ws.Range("I34").Value = 55556732
$ws.Range("K34").Value = 55556732
$ws.Range("M34").Value = -55556530

$ws.Range("H94").Value = 1621.4117
$ws.Range("I94").Value = 187.66667
$ws.Range("J94").Value = 1928.6428
$ws.Range("K94").Value = 187.66667
$ws.Range("L94").Value = 1928.6428
$ws.Range("M94").Value = 263.33333
$ws.Range("N94").Value = -2830.6428

$ws.Range("H99").Value = 12626.934
$ws.Range("I99").Value = 15419.556
$ws.Range("K99").Value = 15419.556
$ws.Range("M99").Value = -13921.556

$ws.Range("H107").Value = 1808.5
$ws.Range("I107").Value = 829.94116
$ws.Range("K107").Value = 829.94116
$ws.Range("M107").Value = 1090.05884

$ws.Range("H122").Value = 2872.1667
$ws.Range("I122").Value = 1105.091
$ws.Range("K122").Value = 3315.273
$ws.Range("M122").Value = -865.2729999999997

$ws.Range("H126").Value = 12626.934
$ws.Range("I126").Value = 15419.556
$ws.Range("K126").Value = 46258.66800000001
$ws.Range("M126").Value = -43788.66800000001

$ws.Range("H132").Value = 43019580
$ws.Range("I132").Value = 55558130
$ws.Range("K132").Value = 166674390
$ws.Range("M132").Value = -166671860

$ws.Range("H134").Value = 1150.85
$ws.Range("I134").Value = 1255.5555
$ws.Range("K134").Value = 3766.6665
$ws.Range("M134").Value = -1231.6665

$ws.Range("H141").Value = 123844.78
$ws.Range("J141").Value = 123844.78
$ws.Range("L141").Value = 123844.78
$ws.Range("N141").Value = -134204.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I2").Value = 65.90000000000001
$ws.Range("J2").Value = 273
$ws.Range("K2").Value = 395.4
$ws.Range("L2").Value = 1638
$ws.Range("M2").Value = -282.4
$ws.Range("N2").Value = -1864

$ws.Range("H8").Value = 95.166664
$ws.Range("I8").Value = 95.166664
$ws.Range("K8").Value = 285.499992
$ws.Range("M8").Value = -146.499992

$ws.Range("H68").Value = 135666.3
$ws.Range("J68").Value = 148278.22
$ws.Range("L68").Value = 444834.66
$ws.Range("N68").Value = -446456.66

$ws.Range("H71").Value = 135666.3
$ws.Range("J71").Value = 148278.22
$ws.Range("L71").Value = 1334503.98
$ws.Range("N71").Value = -1342615.98

$ws.Range("H75").Value = 35719504
$ws.Range("J75").Value = 45460144
$ws.Range("L75").Value = 136380432
$ws.Range("N75").Value = -136382428

$ws.Range("H78").Value = 35719504
$ws.Range("J78").Value = 45460144
$ws.Range("L78").Value = 409141296
$ws.Range("N78").Value = -409151280

$ws.Range("H98").Value = 899.6667
$ws.Range("J98").Value = 999.619
$ws.Range("L98").Value = 2998.857
$ws.Range("N98").Value = -5994.857

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15943586
$ws.Range("I102").Value = 20405238
$ws.Range("K102").Value = 20405238
$ws.Range("M102").Value = -20403616

$ws.Range("H107").Value = 3664137.5
$ws.Range("J107").Value = 782.625
$ws.Range("L107").Value = 782.625
$ws.Range("N107").Value = -4622.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 83337384
$ws.Range("I122").Value = 111115176
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 333345528
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -333343078
$ws.Range("N122").Value = -16900

$ws.Range("H136").Value = 2240.179
$ws.Range("J136").Value = 4809.846
$ws.Range("L136").Value = 14429.538
$ws.Range("N136").Value = -19529.538

$ws.Range("H139").Value = 89234.836
$ws.Range("J139").Value = 97082
$ws.Range("L139").Value = 97082
$ws.Range("N139").Value = -107362

$ws.Range("H140").Value = 91323.09
$ws.Range("J140").Value = 91323.09
$ws.Range("L140").Value = 91323.09
$ws.Range("N140").Value = -101683.09

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2981104.5
$ws.Range("J81").Value = 5000
$ws.Range("L81").Value = 10000
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 2981104.5
$ws.Range("J84").Value = 5000
$ws.Range("L84").Value = 50000
$ws.Range("N84").Value = -60608

$ws.Range("H122").Value = 2788.182
$ws.Range("I122").Value = 2855.1538
$ws.Range("K122").Value = 8565.4614
$ws.Range("M122").Value = -6115.4614

$ws.Range("H132").Value = 19162118
$ws.Range("I132").Value = 2420545.8
$ws.Range("J132").Value = 83338140
$ws.Range("K132").Value = 7261637.399999999
$ws.Range("L132").Value = 250014420
$ws.Range("M132").Value = -7259107.399999999
$ws.Range("N132").Value = -250019480
